$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("STeAM")

# Update the TRUID values embedded in the pipe-delimited payload strings
# in column I (new-user request bodies) for rows 2-5.
$ws.Range("I2").Value = "USER_NAME=Project_NeonUser1@1p.com||USER_PASSWORD=1234qwer$$!||PASSWORD_GENERATE=NO||EMAIL_GENERATE=YES||USER_FIRST_NAME=Project||USER_LAST_NAME=Neon1||USER_MIDDLE_NAME=TR||USER_TRUID=655694d8-f72a-4be4-906c-1e53d3232233"
$ws.Range("I3").Value = "USER_NAME=Project_NeonUser1@1p.com||USER_PASSWORD=1234qwer$$!||PASSWORD_GENERATE=NO||EMAIL_GENERATE=YES||USER_FIRST_NAME=Project||USER_LAST_NAME=Neon1||USER_MIDDLE_NAME=TR||USER_TRUID=655694d8-f72a-4be4-906c-1e53d3230987"
$ws.Range("I5").Value = "USER_NAME=Neon_JDRUser5@1p.com||USER_PASSWORD=1234qwer$$!||PASSWORD_GENERATE=NO||EMAIL_GENERATE=YES||USER_FIRST_NAME=JDR||USER_LAST_NAME=E5||USER_MIDDLE_NAME=REDDY5||USER_TRUID=655694d8-f72a-4be4-906c-1e53d3232233"
$ws.Range("I4").Value = "USER_NAME=Neon_JDRUser6@1p.com||USER_PASSWORD=1234qwer$$!||PASSWORD_GENERATE=NO||EMAIL_GENERATE=YES||USER_TRUID=655694d8-f72a-4be4-906c-1e53d3235566"

# The old workbook had mailto: hyperlinks on I2/I4 pointing at the stale
# request strings; drop them now that the values above are current.
$ws.Hyperlinks.Delete()

# Reflect where the author was last working when the file was saved.
$ws.Range("I4").Select()
